$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New step "PasoGrafica" added to the F column of row 8
$ws.Range("F8").Value = "PasoGrafica"

# Clear the old DW1xx / FUNCION lookup helper table (labels, hex codes and formulas)
$ws.Range("H15:K28").ClearContents()

# Leave selection where the last edit happened
$ws.Range("G15").Select()
